$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set D and E columns (rows 2-51) to text format so numeric-looking
# strings (e.g. "1.003") are preserved as text, matching the source data.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = '28.840.91'
$ws.Cells.Item(2, 5).Value = '  -1.25%  '
$ws.Cells.Item(3, 4).Value = '1.812.68'
$ws.Cells.Item(3, 5).Value = '  -0.90%  '
$ws.Cells.Item(4, 4).Value = '1.003'
$ws.Cells.Item(4, 5).Value = '  +0.24%  '
$ws.Cells.Item(5, 4).Value = '232.49'
$ws.Cells.Item(5, 5).Value = '  -2.04%  '
$ws.Cells.Item(6, 5).Value = '  -3.08%  '
$ws.Cells.Item(7, 5).Value = '  +0.21%  '
$ws.Cells.Item(8, 4).Value = '0.2756'
$ws.Cells.Item(8, 5).Value = '  -2.36%  '
$ws.Cells.Item(9, 4).Value = '0.06741'
$ws.Cells.Item(9, 5).Value = '  -5.03%  '
$ws.Cells.Item(10, 4).Value = '22.87'
$ws.Cells.Item(11, 4).Value = '0.07500'
$ws.Cells.Item(11, 5).Value = '  -1.81%  '
$ws.Cells.Item(12, 4).Value = '1.810.87'
$ws.Cells.Item(12, 5).Value = '  -1.20%  '
$ws.Cells.Item(13, 4).Value = '4.676'
$ws.Cells.Item(13, 5).Value = '  -2.99%  '
$ws.Cells.Item(14, 4).Value = '0.6236'
$ws.Cells.Item(14, 5).Value = '  -2.26%  '
$ws.Cells.Item(15, 4).Value = '0.000009301'
$ws.Cells.Item(15, 5).Value = '  -6.75%  '
$ws.Cells.Item(16, 4).Value = '74.69'
$ws.Cells.Item(16, 5).Value = '  -5.87%  '
$ws.Cells.Item(17, 4).Value = '28.630.92'
$ws.Cells.Item(17, 5).Value = '  -1.93%  '
$ws.Cells.Item(18, 4).Value = '5.451'
$ws.Cells.Item(18, 5).Value = '  -8.53%  '
$ws.Cells.Item(19, 5).Value = '  +0.18%  '
$ws.Cells.Item(20, 4).Value = '208.81'
$ws.Cells.Item(20, 5).Value = '  -8.83%  '
$ws.Cells.Item(21, 4).Value = '11.38'
$ws.Cells.Item(21, 5).Value = '  -3.83%  '
$ws.Cells.Item(22, 4).Value = '6.762'
$ws.Cells.Item(22, 5).Value = '  -4.06%  '
$ws.Cells.Item(23, 5).Value = '  +0.18%  '
$ws.Cells.Item(24, 4).Value = '154.30'
$ws.Cells.Item(24, 5).Value = '  -0.79%  '
$ws.Cells.Item(25, 2).Value = 'Cosmos'
$ws.Cells.Item(25, 3).Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Cells.Item(25, 4).Value = '7.804'
$ws.Cells.Item(25, 5).Value = '  -3.84%  '
$ws.Cells.Item(26, 2).Value = 'Stellar'
$ws.Cells.Item(26, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Cells.Item(26, 4).Value = '0.1271'
$ws.Cells.Item(26, 5).Value = '  -2.58%  '
$ws.Cells.Item(27, 4).Value = '16.33'
$ws.Cells.Item(27, 5).Value = '  -2.30%  '
$ws.Cells.Item(28, 2).Value = 'Toncoin'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Cells.Item(28, 4).Value = '1.406'
$ws.Cells.Item(28, 5).Value = '  -5.32%  '
$ws.Cells.Item(29, 2).Value = 'Hedera'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(29, 4).Value = '0.06285'
$ws.Cells.Item(29, 5).Value = '  -6.93%  '
$ws.Cells.Item(30, 4).Value = '1.428'
$ws.Cells.Item(30, 5).Value = '  -2.14%  '
$ws.Cells.Item(31, 4).Value = '3.736'
$ws.Cells.Item(31, 5).Value = '  -3.10%  '
$ws.Cells.Item(32, 4).Value = '3.695'
$ws.Cells.Item(32, 5).Value = '  -3.75%  '
$ws.Cells.Item(33, 4).Value = '1.697'
$ws.Cells.Item(33, 5).Value = '  -2.59%  '
$ws.Cells.Item(34, 4).Value = '1.051'
$ws.Cells.Item(34, 5).Value = '  -7.14%  '
$ws.Cells.Item(35, 4).Value = '0.6352'
$ws.Cells.Item(35, 5).Value = '  -3.35%  '
$ws.Cells.Item(36, 4).Value = '2.523'
$ws.Cells.Item(36, 5).Value = '  -1.27%  '
$ws.Cells.Item(37, 4).Value = '2.725'
$ws.Cells.Item(37, 5).Value = '  -1.21%  '
$ws.Cells.Item(38, 4).Value = '6.434'
$ws.Cells.Item(38, 5).Value = '  -2.44%  '
$ws.Cells.Item(39, 4).Value = '0.01693'
$ws.Cells.Item(39, 5).Value = '  -4.06%  '
$ws.Cells.Item(40, 4).Value = '1.130.86'
$ws.Cells.Item(40, 5).Value = '  -8.57%  '
$ws.Cells.Item(41, 4).Value = '0.8678'
$ws.Cells.Item(41, 5).Value = '  -6.21%  '
$ws.Cells.Item(42, 5).Value = '  +0.21%  '
$ws.Cells.Item(43, 2).Value = 'Quant'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Cells.Item(43, 4).Value = '100.27'
$ws.Cells.Item(43, 5).Value = '  -0.87%  '
$ws.Cells.Item(44, 2).Value = 'RocketPoolETH'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Cells.Item(44, 4).Value = '1.968.24'
$ws.Cells.Item(44, 5).Value = '  -0.90%  '
$ws.Cells.Item(45, 4).Value = '60.54'
$ws.Cells.Item(45, 5).Value = '  -4.95%  '
$ws.Cells.Item(46, 5).Value = '  -4.61%  '
$ws.Cells.Item(47, 4).Value = '1.571'
$ws.Cells.Item(47, 5).Value = '  -3.51%  '
$ws.Cells.Item(48, 4).Value = '0.4513'
$ws.Cells.Item(48, 5).Value = '  -1.08%  '
$ws.Cells.Item(49, 2).Value = 'EnergySwap'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(49, 4).Value = '8.308'
$ws.Cells.Item(49, 5).Value = '  -3.13%  '
$ws.Cells.Item(50, 2).Value = 'Cronos'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Cells.Item(50, 4).Value = '0.05448'
$ws.Cells.Item(50, 5).Value = '  -1.99%  '
$ws.Cells.Item(51, 5).Value = '  +0.17%  '

# Restore default style on the D:E range (clears the temporary text
# number format applied above) so only cell values differ from the source.
$ws.Range("D2:E51").Style = "Normal"
